$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 2448
$ws.Range("I12").Value = 890
$ws.Range("J12").Value = 2759.6
$ws.Range("K12").Value = 890
$ws.Range("L12").Value = 2759.6
$ws.Range("M12").Value = -720
$ws.Range("N12").Value = -3099.6

$ws.Range("H93").Value = 60000
$ws.Range("J93").Value = 60000
$ws.Range("L93").Value = 60000
$ws.Range("N93").Value = -64992

$ws.Range("H99").Value = 294.30435
$ws.Range("I99").Value = 267.09525
$ws.Range("J99").Value = 580
$ws.Range("K99").Value = 801.28575
$ws.Range("L99").Value = 1740
$ws.Range("M99").Value = 696.71425
$ws.Range("N99").Value = -4736

$ws.Range("H116").Value = 11957
$ws.Range("I116").Value = 5239.8
$ws.Range("J116").Value = 28750
$ws.Range("K116").Value = 5239.8
$ws.Range("L116").Value = 28750
$ws.Range("M116").Value = -1797.8
$ws.Range("N116").Value = -35634

$ws.Range("H132").Value = 4305.7026
$ws.Range("I132").Value = 4590.147
$ws.Range("K132").Value = 13770.441
$ws.Range("M132").Value = -11240.441

$ws.Range("H138").Value = 3478.6667
$ws.Range("I138").Value = 1549
$ws.Range("J138").Value = 3814.261
$ws.Range("K138").Value = 4647
$ws.Range("L138").Value = 11442.783
$ws.Range("M138").Value = 493
$ws.Range("N138").Value = -21722.783

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H54").Value = 30044
$ws.Range("I54").Value = 30044
$ws.Range("K54").Value = 30044
$ws.Range("M54").Value = -29275

$ws.Range("H60").Value = 40000
$ws.Range("I60").Value = 40000
$ws.Range("K60").Value = 40000
$ws.Range("M60").Value = -39267

$ws.Range("H61").Value = 1937.0571
$ws.Range("I61").Value = 1038.9445
$ws.Range("K61").Value = 1038.9445
$ws.Range("M61").Value = -826.9445000000001

$ws.Range("H122").Value = 4418.7383
$ws.Range("I122").Value = 4912.871
$ws.Range("K122").Value = 14738.613
$ws.Range("M122").Value = -12288.613

$ws.Range("H134").Value = 50000
$ws.Range("J134").Value = 50000
$ws.Range("L134").Value = 50000
$ws.Range("N134").Value = -60140

$ws.Range("H136").Value = 1937.0571
$ws.Range("I136").Value = 1038.9445
$ws.Range("K136").Value = 3116.8335
$ws.Range("M136").Value = -566.8335000000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 2999.5
$ws.Range("I22").Value = 2999.5
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 2999.5
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -2826.5
$ws.Range("N22").ClearContents()

$ws.Range("H94").Value = 68970456
$ws.Range("I94").Value = 95243940
$ws.Range("K94").Value = 95243940
$ws.Range("M94").Value = -95243489

$ws.Range("H124").Value = 69899
$ws.Range("J124").Value = 69899
$ws.Range("L124").Value = 69899
$ws.Range("N124").Value = -79719

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1985.5333
$ws.Range("I16").Value = 2021.3334
$ws.Range("J16").Value = 1931.8334
$ws.Range("K16").Value = 2021.3334
$ws.Range("L16").Value = 1931.8334
$ws.Range("M16").Value = -1734.3334
$ws.Range("N16").Value = -2505.8334

$ws.Range("H31").Value = 4498
$ws.Range("I31").Value = 4033.9473
$ws.Range("K31").Value = 4033.9473
$ws.Range("M31").Value = -3738.9473

$ws.Range("H34").Value = 4498
$ws.Range("I34").Value = 4033.9473
$ws.Range("K34").Value = 4033.9473
$ws.Range("M34").Value = -3831.9473

$ws.Range("H43").Value = 38838
$ws.Range("J43").Value = 38838
$ws.Range("L43").Value = 38838
$ws.Range("N43").Value = -39206

$ws.Range("H86").Value = 7730.4546
$ws.Range("I86").Value = 7670.6665
$ws.Range("K86").Value = 7670.6665
$ws.Range("M86").Value = -6547.6665

$ws.Range("H89").Value = 7730.4546
$ws.Range("I89").Value = 7670.6665
$ws.Range("K89").Value = 38353.3325
$ws.Range("M89").Value = -32737.3325

$ws.Range("H101").Value = 38838
$ws.Range("J101").Value = 38838
$ws.Range("L101").Value = 38838
$ws.Range("N101").Value = -45328

$ws.Range("H107").Value = 3572246
$ws.Range("I107").Value = 6250406
$ws.Range("J107").Value = 1365.8334
$ws.Range("K107").Value = 6250406
$ws.Range("L107").Value = 1365.8334
$ws.Range("M107").Value = -6248486
$ws.Range("N107").Value = -5205.8334

$ws.Range("H113").Value = 1985.5333
$ws.Range("I113").Value = 2021.3334
$ws.Range("J113").Value = 1931.8334
$ws.Range("K113").Value = 2021.3334
$ws.Range("L113").Value = 1931.8334
$ws.Range("M113").Value = 148.6666
$ws.Range("N113").Value = -6271.8334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1355.4073
$ws.Range("I2").Value = 271
$ws.Range("J2").Value = 2100.9375
$ws.Range("K2").Value = 1626
$ws.Range("L2").Value = 12605.625
$ws.Range("M2").Value = -1513
$ws.Range("N2").Value = -12831.625

$ws.Range("H60").Value = 1113865.4
$ws.Range("I60").Value = 2500736
$ws.Range("J60").Value = 4368.8
$ws.Range("K60").Value = 7502208
$ws.Range("L60").Value = 13106.4
$ws.Range("M60").Value = -7501957
$ws.Range("N60").Value = -13608.4

$ws.Range("H98").Value = 707.5454999999999
$ws.Range("J98").Value = 692.1111
$ws.Range("L98").Value = 2076.3333
$ws.Range("N98").Value = -5072.3333

$ws.Range("H133").Value = 4575
$ws.Range("I133").Value = 2025
$ws.Range("K133").Value = 6075
$ws.Range("M133").Value = -1015

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 50000000
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

$ws.Range("H22").Value = 1455.8
$ws.Range("I22").Value = 426.66666
$ws.Range("J22").Value = 2999.5
$ws.Range("K22").Value = 426.66666
$ws.Range("L22").Value = 2999.5
$ws.Range("M22").Value = 102.33334
$ws.Range("N22").Value = -4057.5

$ws.Range("H81").Value = 50000000
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 50000000
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H41").Value = 30033
$ws.Range("I41").Value = 30033
$ws.Range("K41").Value = 30033
$ws.Range("M41").Value = -29595

$ws.Range("H100").Value = 5156.0835
$ws.Range("J100").Value = 7247
$ws.Range("L100").Value = 7247
$ws.Range("N100").Value = -8329

$ws.Range("H104").Value = 43058
$ws.Range("J104").Value = 43058
$ws.Range("L104").Value = 43058
$ws.Range("N104").Value = -50046

$ws.Range("H132").Value = 6068.625
$ws.Range("I132").Value = 2246.2144
$ws.Range("K132").Value = 6738.6432
$ws.Range("M132").Value = -4208.6432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 7000
$ws.Range("J47").Value = 7000
$ws.Range("L47").Value = 7000
$ws.Range("N47").Value = -8144

$ws.Range("H92").Value = 106599.6
$ws.Range("J92").Value = 106599.6
$ws.Range("L92").Value = 106599.6
$ws.Range("N92").Value = -111591.6

$ws.Range("H100").Value = 142857810
$ws.Range("I100").Value = 835
$ws.Range("J100").Value = 250000530
$ws.Range("K100").Value = 1670
$ws.Range("L100").Value = 500001060
$ws.Range("M100").Value = -1129
$ws.Range("N100").Value = -500002142

$ws.Range("H132").Value = 2369.075
$ws.Range("I132").Value = 2570.8708
$ws.Range("K132").Value = 7712.6124
$ws.Range("M132").Value = -5182.6124
